# Apply the changes described by the diff:
# 1. Change the "boson" column (E2:E12) text values from "W+" to "W"
# 2. Turn on iterative calculation with a max change (delta) of 1E-4
#    (mirrors calcPr iterateDelta="1E-4" in the target workbook)
# 3. Move the active selection to F19

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update boson values in column E (rows 2-12) from "W+" to "W"
$ws.Range("E2:E12").Value = "W"

# 2. Enable iterative calculation with a max change of 1E-4
$excel.Iteration = $true
$excel.MaxIterations = 100
$excel.MaxChange = 0.0001

# 3. Set the active selection to F19
$ws.Range("F19").Select()
